# "added .net core release"
# Slide 32 ("Windows Desktop") - append a new paragraph to the
# "Content Placeholder 2" body listing the current Entity Framework 6
# release timing, split across three runs (middle run carries no
# "dirty" proofing flag in the source deck).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(32)

$shp = $null
foreach ($candidate in $s.Shapes) {
    if ($candidate.Name -eq "Content Placeholder 2") {
        $shp = $candidate
        break
    }
}

$tr = $shp.TextFrame.TextRange
$priorLength = $tr.Length

$part1 = "Release (as of "
$part2 = "this writing) "
$part3 = "September 2019"

# Start a new paragraph after the existing "Entity Framework 6" line,
# then populate it with the full sentence in one shot.
$inserted = $tr.InsertAfter("`r" + $part1 + $part2 + $part3)

# Re-slice the freshly inserted paragraph into three separate runs so
# each piece can carry its own run properties, matching how the text
# was built up in separate edits.
$newParaStart = $priorLength + 2

$run1 = $tr.Characters($newParaStart, $part1.Length)
$run1.Text = $part1

$run2 = $tr.Characters($newParaStart + $part1.Length, $part2.Length)
$run2.Text = $part2

$run3 = $tr.Characters($newParaStart + $part1.Length + $part2.Length, $part3.Length)
$run3.Text = $part3

Write-Output "Paragraphs now: $($tr.Paragraphs().Count)"
Write-Output "New paragraph text: $($tr.Paragraphs(4,1).Text)"
